$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Novos resultados da QUINA (concursos 6908 a 6926), adicionados apos a ultima linha existente (416).
$novosConcursos = @(
  @(6908, 1, 10, 42, 58, 72),
  @(6909, 9, 25, 57, 63, 75),
  @(6910, 2, 18, 20, 49, 68),
  @(6911, 3, 51, 56, 59, 72),
  @(6912, 4, 43, 44, 58, 62),
  @(6913, 25, 31, 38, 42, 58),
  @(6914, 25, 41, 48, 49, 66),
  @(6915, 5, 19, 21, 51, 66),
  @(6916, 8, 54, 58, 72, 76),
  @(6917, 11, 15, 29, 48, 57),
  @(6918, 9, 21, 24, 63, 69),
  @(6919, 4, 6, 9, 26, 64),
  @(6920, 4, 28, 34, 42, 47),
  @(6921, 26, 50, 69, 74, 77),
  @(6922, 16, 26, 36, 51, 56),
  @(6923, 18, 34, 41, 57, 63),
  @(6924, 4, 13, 49, 52, 66),
  @(6925, 9, 25, 44, 46, 62),
  @(6926, 14, 29, 40, 79, 80)
)

$startRow = 417
for ($i = 0; $i -lt $novosConcursos.Count; $i++) {
  $row = $startRow + $i
  $vals = $novosConcursos[$i]
  for ($c = 0; $c -lt $vals.Count; $c++) {
    $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
  }
}

$lastRow = $startRow + $novosConcursos.Count - 1
$selRange = "B" + $startRow + ":F" + $lastRow
$ws.Range($selRange).Select()
